$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44330
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 861

# Row 5
$ws.Range("D5").Value = 44707
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 12500
$ws.Range("Q5").Value = "$/caja 12 kilos empedrada"
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 1042
$ws.Range("T5").Value = 12

# Row 6
$ws.Range("D6").Value = 44334
$ws.Range("N6").Value = 11000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 11500
$ws.Range("Q6").Value = "$/caja 12 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 11500
$ws.Range("T6").Value = 1

# Row 7
$ws.Range("D7").Value = 44719
$ws.Range("M7").Value = 50
$ws.Range("P7").Value = 14400
$ws.Range("R7").Value = "Región del Maule"
$ws.Range("S7").Value = 800

# Row 8
$ws.Range("D8").Value = 44714
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("Q8").Value = "$/caja 18 kilos granel"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 806
$ws.Range("T8").Value = 18
